# Update cryptos list — GitHub Actions scheduled refresh.
#
# All the source cells in this sheet are plain text (inlineStr) even when
# their content looks like a number ("593.99", "68.214.70", ...). A bare
# `Range.Value = "..."` assignment lets Excel's type-sniffer turn anything
# that parses as a number into a real numeric cell, which would change the
# cell's stored type/shape versus the original workbook. Force the cell to
# Text before writing, then hand the style back to "Normal" so we don't
# leave a stray custom number format sitting on the cell afterwards.
function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: Price (column D) and Volume(1h) (column E).
# $null means that column is left untouched for that row.
$updates = @(
    [PSCustomObject]@{ Row=2;  D="68.214.70"; E="  +2.54%  " }
    [PSCustomObject]@{ Row=3;  D="2.531.28";  E="  +0.74%  " }
    [PSCustomObject]@{ Row=4;  D=$null;       E="  +0.07%  " }
    [PSCustomObject]@{ Row=5;  D="593.99";    E="  +1.92%  " }
    [PSCustomObject]@{ Row=6;  D="178.86";    E="  +4.57%  " }
    [PSCustomObject]@{ Row=7;  D=$null;       E="  +0.03%  " }
    [PSCustomObject]@{ Row=8;  D=$null;       E="  +1.23%  " }
    [PSCustomObject]@{ Row=9;  D="2.530.83";  E="  +0.76%  " }
    [PSCustomObject]@{ Row=10; D=$null;       E="  +2.52%  " }
    [PSCustomObject]@{ Row=11; D=$null;       E="  +2.72%  " }
    [PSCustomObject]@{ Row=12; D="5.12";      E="  -0.02%  " }
    [PSCustomObject]@{ Row=13; D=$null;       E="  -0.43%  " }
    [PSCustomObject]@{ Row=14; D="26.88";     E="  +0.61%  " }
    [PSCustomObject]@{ Row=15; D="2.996.48";  E="  +1.33%  " }
    [PSCustomObject]@{ Row=16; D=$null;       E="  +2.15%  " }
    [PSCustomObject]@{ Row=17; D="68.225.62"; E="  +2.77%  " }
    [PSCustomObject]@{ Row=18; D="2.517.88";  E="  +0.34%  " }
    [PSCustomObject]@{ Row=19; D="8.00";      E="  +2.44%  " }
    [PSCustomObject]@{ Row=20; D="11.54";     E="  +2.43%  " }
    [PSCustomObject]@{ Row=21; D="368.51";    E="  +5.93%  " }
    [PSCustomObject]@{ Row=22; D=$null;       E="  +0.25%  " }
    [PSCustomObject]@{ Row=23; D=$null;       E="  +1.65%  " }
    [PSCustomObject]@{ Row=24; D=$null;       E="  -1.41%  " }
    [PSCustomObject]@{ Row=25; D=$null;       E="  -0.03%  " }
    [PSCustomObject]@{ Row=26; D="70.72";     E="  +1.21%  " }
    [PSCustomObject]@{ Row=27; D=$null;       E="  +3.14%  " }
    [PSCustomObject]@{ Row=28; D="2.661.49";  E=$null }
    [PSCustomObject]@{ Row=29; D=$null;       E="  -0.37%  " }
    [PSCustomObject]@{ Row=30; D="0.0₃0997";  E="  +2.10%  " }
    [PSCustomObject]@{ Row=31; D="541.48";    E="  +3.32%  " }
    [PSCustomObject]@{ Row=32; D=$null;       E="  +2.68%  " }
    [PSCustomObject]@{ Row=33; D=$null;       E="  +1.91%  " }
    [PSCustomObject]@{ Row=34; D="1.87";      E="  +2.02%  " }
    [PSCustomObject]@{ Row=35; D=$null;       E="  -0.79%  " }
    [PSCustomObject]@{ Row=36; D=$null;       E="  +0.05%  " }
    [PSCustomObject]@{ Row=37; D=$null;       E="  +0.02%  " }
    [PSCustomObject]@{ Row=38; D="157.50";    E="  +0.24%  " }
    [PSCustomObject]@{ Row=39; D="18.87";     E="  +1.26%  " }
    [PSCustomObject]@{ Row=40; D="18.70";     E="  +1.73%  " }
    [PSCustomObject]@{ Row=41; D="0.356";     E="  +0.33%  " }
    [PSCustomObject]@{ Row=42; D=$null;       E="  +1.06%  " }
    [PSCustomObject]@{ Row=43; D="5.22";      E="  +2.75%  " }
    [PSCustomObject]@{ Row=44; D="2.57";      E="  +3.09%  " }
    [PSCustomObject]@{ Row=45; D=$null;       E="  -0.08%  " }
    [PSCustomObject]@{ Row=46; D="147.55";    E="  -0.82%  " }
    [PSCustomObject]@{ Row=47; D="0.562";     E="  +0.79%  " }
    [PSCustomObject]@{ Row=50; D=$null;       E="  -0.56%  " }
    [PSCustomObject]@{ Row=51; D=$null;       E="  +0.56%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue $ws "D$($u.Row)" $u.D
    }
    if ($null -ne $u.E) {
        Set-TextValue $ws "E$($u.Row)" $u.E
    }
}

# Rows 48 and 49 swap ranking order: BabyDogeCoin drops below Filecoin.
Set-TextValue $ws "B48" "Filecoin"
Set-TextValue $ws "C48" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D48" "3.73"
Set-TextValue $ws "E48" "  +1.20%  "

Set-TextValue $ws "B49" "BabyDogeCoin"
Set-TextValue $ws "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D49" "0.0₆0278"
Set-TextValue $ws "E49" "  +3.21%  "
